$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period text on row 2 (merged A2:G2 - only the top-left anchor
# cell of a merged range actually holds/accepts the value)
$ws.Range("A2").Value = "Период: 2023-11-01 - 2023-11-30"

# Update the remaining data row (row 8) with new values
# A8/D8 look like plain numbers - use a leading apostrophe to force text
# (matching the original file, which stores them as text), then restore
# the original cell formatting via a format-only paste.
$ws.Range("A8").Value = "'70004"
$ws.Range("C8").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B8").Value = "Сбербанк"
$ws.Range("C8").Value = "лебенков"

$ws.Range("D8").Value = "'11300"
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("E8").Value = "01.11.2023 00:00:00"
$ws.Range("F8").Value = "Поступление"
$ws.Range("G8").Value = "Докторская, Волковыское"

# Remove rows 9 and 10 (old extra data rows)
$ws.Range("A9:G10").EntireRow.Delete()

# Widen column G (closest achievable width to the target 25.1897583007813 -
# Excel snaps ColumnWidth to whole-pixel steps, so this is the nearest grid value)
$ws.Range("G1").ColumnWidth = 24.3
